$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "data": append a new row 3 (same shape as row 2) with the new
# record (Xmas 2022 / Xmas vacation 2022).
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("data")

# Copy A2's formatting (bold header-style font/border/alignment) onto A3
# so the new ID cell matches the styling already used for A2, instead of
# allocating a brand-new (duplicate) style entry.
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(3, 1).PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(3, 1).Value = "b979ba9a-87ca-4b10-9b53-48dd201ed68a"

# "01/02/2023" must stay literal text (like B2's "12/10/2022"), not get
# reinterpreted as a date serial number - force it in as text via a
# leading apostrophe, then drop the resulting quote-prefix formatting so
# the cell ends up with the default (unstyled) look, same as B2.
$ws.Cells.Item(3, 2).Value = "'01/02/2023"
$ws.Cells.Item(3, 2).ClearFormats()

# C3 mirrors C2: present but empty.
$ws.Cells.Item(3, 3).NumberFormat = "@"
$ws.Cells.Item(3, 3).Value = ""
$ws.Cells.Item(3, 3).ClearFormats()

$ws.Cells.Item(3, 4).Value = "Xmas 2022"
$ws.Cells.Item(3, 5).Value = "Xmas 2022"
$ws.Cells.Item(3, 6).Value = "Xmas vacation 2022"

# ---------------------------------------------------------------------
# Sheet "headers": row 7 (description) - D/E/F move from the literal
# text "TRUE" to real boolean TRUE values.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("headers")
$ws2.Cells.Item(7, 4).Value = $true
$ws2.Cells.Item(7, 5).Value = $true
$ws2.Cells.Item(7, 6).Value = $true
